$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.404.32"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "1.711.40"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.32"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5326"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06610"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.93"
$ws.Range("E10").Value = "  -4.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07644"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.567"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").Value = "1.702.73"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").Value = "1.948.93"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5775"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "0.0₅8190"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "27.397.74"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.49"
$ws.Range("E19").Value = "  -3.69%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.673"
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.981"
$ws.Range("E23").Value = "  -3.91%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.52"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.739"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1218"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.297"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("E29").Value = "  -4.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05416"
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.292"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.513"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.433"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.649"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.876"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9506"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.417"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5873"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01635"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.869"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "1.046.39"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8418"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.95"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Value = "1.854.70"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.04"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4512"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.076"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("E51").Value = "  -1.45%  "
